$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A27:E27").Copy()
$ws.Range("A28:E28").PasteSpecial(-4122)

$ws.Cells.Item(28, 1).Value = 26
$ws.Cells.Item(28, 2).Value = "CounterMagic"
$ws.Cells.Item(28, 3).Value = "PassiveSkill"
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = 0

$ws.Rows.Item(28).Select()
